# Restore C10 on the "Rules" sheet to its updated value (18 -> 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
